$d = $word.ActiveDocument

$pairs = @(
    @("2025-10-31 Friday", "2025-11-01 Saturday"),
    @("787÷3=", "691÷3="),
    @("307÷5=", "136÷9="),
    @("633÷5=", "270÷9="),
    @("169÷3=", "425÷6="),
    @("136÷6=", "988÷8="),
    @("830÷2=", "468÷6="),
    @("121÷3=", "781÷9="),
    @("303÷3=", "437÷5="),
    @("580÷4=", "655÷5="),
    @("990÷6=", "968÷5="),
    @("124÷6=", "275÷5="),
    @("682÷6=", "781÷8="),
    @("490÷3=", "482÷4="),
    @("987÷5=", "407÷9="),
    @("457÷3=", "781÷5="),
    @("990÷3=", "876÷9="),
    @("824÷2=", "718÷8="),
    @("695÷8=", "957÷2="),
    @("139÷3=", "155÷2="),
    @("298÷2=", "793÷3="),
    @("168÷8=", "292÷6="),
    @("350÷9=", "772÷2="),
    @("603÷3=", "362÷7="),
    @("829÷4=", "295÷6="),
    @("295÷7=", "117÷5=")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
